# The "Saldo" export sheet lost two account rows from the report:
#   005232019 - PEDRO     - 3000
#   004853111 - MARCONDES - 2459.89
# Locate them by their account number (column A) so the edit is robust
# even if row positions shift, then delete the whole rows so everything
# below moves up (no blank rows left behind).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$accountsToRemove = @("005232019", "004853111")

$rowsToDelete = @()
foreach ($account in $accountsToRemove) {
    $hit = $ws.Columns("A").Find($account)
    if ($hit -ne $null) {
        $rowsToDelete += $hit.Row
    }
}

# Delete from the bottom-most row upward so earlier row numbers stay valid.
$rowsToDelete = $rowsToDelete | Sort-Object -Descending

foreach ($r in $rowsToDelete) {
    $ws.Rows($r).Delete()
}
